# Update the EC database: reorder rows 16-27 by "Periodo Mora" (column E)
# ascending instead of descending, keeping each period's "Valor Mora" (col F)
# and other row data tied to its own period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# Valor Mora (col F): period 2401 keeps 19067, all the others keep 52000.
$ws.Cells.Item(16, 6).Value = 52000
$ws.Cells.Item(27, 6).Value = 19067
